# Fruta / hortaliza, semanal
# Insert a new weekly record above the current row 133 (pushing the
# existing rows 133-169 down to 134-170) on the single data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 133; Excel shifts rows 133:169 down
# to 134:170 and the sheet's used range grows to A1:T170.
$ws.Rows.Item(133).Insert()

# Populate the newly inserted row 133 with the new weekly price record.
$ws.Cells.Item(133, 1).Value  = 5
$ws.Cells.Item(133, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(133, 3).Value  = "Maule"
$ws.Cells.Item(133, 4).Value  = 44508
$ws.Cells.Item(133, 5).Value  = 7
$ws.Cells.Item(133, 6).Value  = "Fruta"
$ws.Cells.Item(133, 7).Value  = 100108
$ws.Cells.Item(133, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(133, 9).Value  = 100108005
$ws.Cells.Item(133, 10).Value = "Piña"
$ws.Cells.Item(133, 11).Value = "Caramelo"
$ws.Cells.Item(133, 12).Value = "Segunda"
$ws.Cells.Item(133, 13).Value = 540
$ws.Cells.Item(133, 14).Value = 17000
$ws.Cells.Item(133, 15).Value = 17000
$ws.Cells.Item(133, 16).Value = 17000
$ws.Cells.Item(133, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(133, 18).Value = "Ecuador"
$ws.Cells.Item(133, 19).Value = 1214
$ws.Cells.Item(133, 20).Value = 14
